# Populate the first data row (row 2) of the "Property" sheet with the
# ProxyServer's configuration data, then move the active selection to G4,
# matching the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values in the order that reproduces the expected shared-string
# insertion order: 127.0.0.1, ProxyServer_1, 000105001.
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "ProxyServer_1"
$ws.Range("B2").Value = "000105001"
$ws.Range("C2").Value = "ProxyServer_1"
# Give the newly-populated C2 the same (text) number format as A2/B2 so it
# reuses the existing "text" cell style instead of the default style.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 5001

# Move the selection/active cell to G4, as in the edited workbook.
$ws.Range("G4").Select()
